{"js": "// The author fixed a typo in the closing paragraph (\"the plan sing\" ->\n// \"the plane sing\" \u2014 a missing \"e\" in \"plane\") and cleaned up the stray\n// \"_GoBack\" bookmark that Word leaves behind after an editing session.\n\nconst body = context.document.body;\n\n// 1) Fix the typo: \"the plan sing\" -> \"the plane sing\".\nconst hits = body.search(\"the plan sing\", { matchCase: true });\nhits.load(\"items,text\");\nawait context.sync();\n\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"the plane sing\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Remove the leftover \"_GoBack\" bookmark from the last paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The author fixed a typo in the closing paragraph (\"the plan sing\" ->\n# \"the plane sing\" \u2014 a missing \"e\" in \"plane\") and cleaned up the stray\n# \"_GoBack\" bookmark that Word leaves behind after an editing session.\n\n$d = $word.ActiveDocument\n\n# 1) Fix the typo: \"the plan sing\" -> \"the plane sing\".\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = \"the plan sing\"\n$find.MatchCase = $true\n$found = $find.Execute()\nif ($found) {\n  $rng.Text = \"the plane sing\"\n}\n\n# 2) Remove the leftover \"_GoBack\" bookmark from the last paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
